# Task: nvidia-stock-analysis
# Change the time duration of the needed data on the "Basic Trend" sheet:
# replace the quarter list 2023Q1..2024Q4 (8 quarters) with 2024Q3..2025Q2
# (4 quarters), shrinking the used range from A1:D9 down to A1:D5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Basic Trend")

# Overwrite the first four quarter rows with the new quarter labels.
$ws.Range("A2").Value = "2024Q3"
$ws.Range("A3").Value = "2024Q4"
$ws.Range("A4").Value = "2025Q1"
$ws.Range("A5").Value = "2025Q2"

# Drop the now-unused trailing quarter rows (2024Q1..2024Q4 previously in
# rows 6-9) so the sheet only spans A1:D5.
$ws.Range("A6:A9").ClearContents() | Out-Null

# Leave the selection where the author left it after editing.
$ws.Range("B11").Select() | Out-Null
